$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- "Auto eval" row (row 18) -------------------------------------------
# A18: classification "Analyse" (re-uses the existing shared string)
$ws.Range("A18").Value = "Analyse"

# C18: duration in minutes for the self-evaluation entry
$ws.Range("C18").Value = 25

# D18: explanation text for the new entry (new shared string)
$ws.Range("D18").Value = "Feet back du livrable de vendredi dernier"

# C66 (total hours) is a live formula -> MROUND(SUM(C6:C65)/60,0.2); it will
# automatically recalculate to reflect the new C18 value once the engine
# recalcs after this script runs.

# --- Selection / window bookkeeping --------------------------------------
# Move the active selection in the frozen "Journal" view to D19, matching
# the saved cursor position recorded in the workbook.
[void]$ws.Range("D19").Select()

# Restore a "Normal" (non page-break-preview) view for the active window,
# matching the author's saved view state.
$excel.ActiveWindow.View = 1

# Window geometry recorded by Excel on save.
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = -120
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15840
